$wb = $excel.ActiveWorkbook

$ws_Estimated = $wb.Worksheets.Item("Estimated")
$ws_Estimated.Range("B2").Value = "Fri Mar 08 22:31:55 EST 2024"
$ws_Estimated.Range("B3").Value = "Fri Mar 08 22:32:35 EST 2024"
$ws_Estimated.Range("B4").Value = "Fri Mar 08 22:33:11 EST 2024"
$ws_Estimated.Range("B5").Value = "Fri Mar 08 22:33:46 EST 2024"
$ws_Estimated.Range("B6").Value = "Fri Mar 08 22:34:21 EST 2024"
$ws_Estimated.Range("B7").Value = "Fri Mar 08 22:34:57 EST 2024"
$ws_Existing = $wb.Worksheets.Item("Existing")
$ws_Existing.Range("B2").Value = "Fri Mar 08 22:35:32 EST 2024"
$ws_Existing.Range("B3").Value = "Fri Mar 08 22:36:08 EST 2024"
$ws_Existing.Range("B4").Value = "Fri Mar 08 22:36:43 EST 2024"
$ws_Existing.Range("B5").Value = "Fri Mar 08 22:37:19 EST 2024"
$ws_Existing.Range("B6").Value = "Fri Mar 08 22:37:55 EST 2024"
$ws_Existing.Range("B7").Value = "Fri Mar 08 22:38:30 EST 2024"
$ws_Existing.Range("B8").Value = "Fri Mar 08 22:39:06 EST 2024"
$ws_Existing.Range("B9").Value = "Fri Mar 08 22:39:43 EST 2024"
$ws_Existing.Range("B10").Value = "Fri Mar 08 22:40:20 EST 2024"
$ws_Existing.Range("B11").Value = "Fri Mar 08 22:40:57 EST 2024"
$ws_Existing.Range("B12").Value = "Fri Mar 08 22:41:32 EST 2024"
$ws_Existing.Range("B13").Value = "Fri Mar 08 22:42:08 EST 2024"
$ws_Existing.Range("B14").Value = "Fri Mar 08 22:42:43 EST 2024"
$ws_Existing.Range("B15").Value = "Fri Mar 08 22:43:19 EST 2024"
$ws_Existing.Range("B16").Value = "Fri Mar 08 22:43:54 EST 2024"
$ws_Existing.Range("B17").Value = "Fri Mar 08 22:44:30 EST 2024"
$ws_Existing.Range("B18").Value = "Fri Mar 08 22:45:05 EST 2024"
$ws_Existing.Range("B19").Value = "Fri Mar 08 22:45:41 EST 2024"
$ws_NewTaxReturn = $wb.Worksheets.Item("NewTaxReturn")
$ws_NewTaxReturn.Range("B2").Value = "Fri Mar 08 22:46:17 EST 2024"
$ws_NewTaxReturn.Range("B3").Value = "Fri Mar 08 22:46:52 EST 2024"
$ws_NewTaxReturn.Range("B4").Value = "Fri Mar 08 22:47:27 EST 2024"
$ws_NewTaxReturn.Range("B5").Value = "Fri Mar 08 22:48:02 EST 2024"
$ws_NewTaxReturn.Range("B6").Value = "Fri Mar 08 22:48:37 EST 2024"
$ws_NewTaxReturn.Range("B7").Value = "Fri Mar 08 22:49:13 EST 2024"
$ws_NewTaxReturn.Range("B8").Value = "Fri Mar 08 22:49:48 EST 2024"
$ws_NewTaxReturn.Range("B9").Value = "Fri Mar 08 22:50:23 EST 2024"
$ws_NewTaxReturn.Range("B10").Value = "Fri Mar 08 22:50:57 EST 2024"
$ws_NewTaxReturn.Range("B11").Value = "Fri Mar 08 22:51:32 EST 2024"
$ws_NewTaxReturn.Range("B12").Value = "Fri Mar 08 22:52:07 EST 2024"
$ws_NewTaxReturn.Range("B13").Value = "Fri Mar 08 22:52:42 EST 2024"
$ws_NewTaxReturn.Range("B14").Value = "Fri Mar 08 22:53:16 EST 2024"
$ws_NewTaxReturn.Range("B15").Value = "Fri Mar 08 22:53:51 EST 2024"
$ws_NewTaxReturn.Range("B16").Value = "Fri Mar 08 22:54:26 EST 2024"
$ws_NewTaxReturn.Range("B17").Value = "Fri Mar 08 22:55:01 EST 2024"
$ws_NewTaxReturn.Range("B18").Value = "Fri Mar 08 22:55:38 EST 2024"
$ws_NewTaxReturn.Range("B19").Value = "Fri Mar 08 22:56:15 EST 2024"
$ws_NewTaxReturn.Range("B20").Value = "Fri Mar 08 22:56:52 EST 2024"
$ws_NewTaxReturn.Range("B21").Value = "Fri Mar 08 22:57:27 EST 2024"
$ws_NewTaxReturn.Range("B22").Value = "Fri Mar 08 22:58:04 EST 2024"
$ws_NewTaxReturn.Range("B23").Value = "Fri Mar 08 22:58:41 EST 2024"
$ws_NewTaxReturn.Range("B24").Value = "Fri Mar 08 22:59:16 EST 2024"
$ws_NewTaxReturn.Range("B25").Value = "Fri Mar 08 22:59:53 EST 2024"
$ws_NewTaxReturn.Range("B26").Value = "Fri Mar 08 23:00:31 EST 2024"
$ws_NewTaxReturn.Range("B27").Value = "Fri Mar 08 23:01:08 EST 2024"
$ws_NewTaxReturn.Range("B28").Value = "Fri Mar 08 23:01:45 EST 2024"
$ws_NewTaxReturn.Range("B29").Value = "Fri Mar 08 23:02:22 EST 2024"
$ws_NewTaxReturn.Range("B30").Value = "Fri Mar 08 23:02:59 EST 2024"
$ws_NewTaxReturn.Range("B31").Value = "Fri Mar 08 23:03:36 EST 2024"
$ws_NewTaxReturn.Range("B32").Value = "Fri Mar 08 23:04:14 EST 2024"
$ws_NewTaxReturn.Range("B33").Value = "Fri Mar 08 23:04:49 EST 2024"
$ws_NewTaxReturn.Range("B34").Value = "Fri Mar 08 23:05:26 EST 2024"
$ws_NewTaxReturn.Range("B35").Value = "Fri Mar 08 23:06:03 EST 2024"
$ws_NewTaxReturn.Range("B36").Value = "Fri Mar 08 23:06:39 EST 2024"
$ws_NewTaxReturn.Range("B37").Value = "Fri Mar 08 23:07:16 EST 2024"
$ws_NewTaxReturn.Range("B38").Value = "Fri Mar 08 23:07:54 EST 2024"
$ws_NewTaxReturn.Range("B39").Value = "Fri Mar 08 23:08:31 EST 2024"
$ws_NewTaxReturn.Range("B40").Value = "Fri Mar 08 23:09:08 EST 2024"
$ws_NewTaxReturn.Range("B41").Value = "Fri Mar 08 23:09:45 EST 2024"
$ws_NewTaxReturn.Range("B42").Value = "Fri Mar 08 23:10:23 EST 2024"
$ws_NewTaxReturn.Range("B43").Value = "Fri Mar 08 23:11:01 EST 2024"
$ws_NewTaxReturn.Range("B44").Value = "Fri Mar 08 23:11:38 EST 2024"
$ws_NewTaxReturn.Range("B45").Value = "Fri Mar 08 23:12:13 EST 2024"
$ws_NewTaxReturn.Range("B46").Value = "Fri Mar 08 23:12:50 EST 2024"
$ws_NewTaxReturn.Range("B47").Value = "Fri Mar 08 23:13:27 EST 2024"
$ws_NewTaxReturn.Range("B48").Value = "Fri Mar 08 23:14:03 EST 2024"
$ws_NewTaxReturn.Range("B49").Value = "Fri Mar 08 23:14:41 EST 2024"
$ws_NewTaxReturn.Range("B50").Value = "Fri Mar 08 23:15:18 EST 2024"
$ws_NewTaxReturn.Range("B51").Value = "Fri Mar 08 23:15:55 EST 2024"
$ws_NewTaxReturn.Range("B52").Value = "Fri Mar 08 23:16:32 EST 2024"
$ws_Personal_IND = $wb.Worksheets.Item("Personal_IND")
$ws_Personal_IND.Range("B2").Value = "Fri Mar 08 23:17:10 EST 2024"
$ws_Personal_IND.Range("B4").Value = "Fri Mar 08 23:17:43 EST 2024"
$ws_Personal_IND.Range("B5").Value = "Fri Mar 08 23:18:17 EST 2024"
$ws_Personal_IND.Range("B6").Value = "Fri Mar 08 23:18:48 EST 2024"
$ws_Personal_IND.Range("B7").Value = "Fri Mar 08 23:19:19 EST 2024"
$ws_Personal_IND.Range("B8").Value = "Fri Mar 08 23:19:50 EST 2024"
$ws_Personal_IND.Range("B9").Value = "Fri Mar 08 23:20:22 EST 2024"
$ws_Personal_JNT = $wb.Worksheets.Item("Personal_JNT")
$ws_Personal_JNT.Range("B2").Value = "Fri Mar 08 23:20:53 EST 2024"
$ws_Personal_JNT.Range("B4").Value = "Fri Mar 08 23:21:25 EST 2024"
$ws_Personal_JNT.Range("B5").Value = "Fri Mar 08 23:21:56 EST 2024"
$ws_Personal_JNT.Range("B6").Value = "Fri Mar 08 23:22:35 EST 2024"
$ws_Personal_EL = $wb.Worksheets.Item("Personal_EL")
$ws_Personal_EL.Range("B2").Value = "Fri Mar 08 23:23:13 EST 2024"
$ws_Personal_EL.Range("B3").Value = "Fri Mar 08 23:23:52 EST 2024"
